$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing asterisk markers from these cells (renumbering the footnotes)
$ws.Range("D34").Value = "*Refunded `$23.58 after 2 motors were out of stock"
$ws.Range("D35").Value = "**Refunded `$8.18 for an incorrect shipment, an additional `$326.00 for returned USB cables, and another `$999.00 for returned power adapters"
$ws.Range("D20").Value = "Hobby King Order 4.pdf*"
$ws.Range("D22").Value = "Amazon Order 9.pdf**"
$ws.Range("A11").Value = "Amazon Order 6"

# Remove the "not sure whether to expense" note entirely
$ws.Range("A35").ClearContents()

# Restore the selection/scroll state to A12 (no frozen top-left override)
$ws.Range("A12").Select()
